# Weekly refresh of the Granada price table: the newest week's records
# are inserted at the top (rows 2-4) and all of last week's rows shift
# down by three positions (old row N -> row N+3), extending the table
# from 16 data rows (2-17) to 19 data rows (2-20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the sheet by 3 rows. Inserting right after the existing data
# (rather than at the top) means the new blank rows simply inherit the
# plain formatting of the row above them instead of the bold header
# style, so no extra styles are introduced.
$ws.Range("A18:T20").Insert(-4121)

# Re-write the full data block (rows 2-20) in its final order: the new
# week first, followed by the previously existing weeks shifted down.
$final = @(
    @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44664, 13, "Fruta", 100104, "Frutos de pepita", 100104001, "Granada", "Sin especificar", "Especial", 300, 21600, 21600, 21600, "`$/caja 18 kilos granel", "Provincia de Limarí", 1200, 18),
    @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44664, 13, "Fruta", 100104, "Frutos de pepita", 100104001, "Granada", "Sin especificar", "Primera", 250, 18000, 18000, 18000, "`$/caja 18 kilos granel", "Provincia de Limarí", 1000, 18),
    @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44664, 13, "Fruta", 100104, "Frutos de pepita", 100104001, "Granada", "Sin especificar", "Segunda", 250, 16000, 16000, 16000, "`$/caja 18 kilos granel", "Provincia de Limarí", 889, 18),
    @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44658, 13, "Fruta", 100104, "Frutos de pepita", 100104001, "Granada", "Sin especificar", "Especial", 280, 21600, 21600, 21600, "`$/caja 18 kilos granel", "Provincia de Limarí", 1200, 18),
    @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44658, 13, "Fruta", 100104, "Frutos de pepita", 100104001, "Granada", "Sin especificar", "Primera", 330, 16200, 16200, 16200, "`$/caja 18 kilos granel", "Provincia de Limarí", 900, 18),
    @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44658, 13, "Fruta", 100104, "Frutos de pepita", 100104001, "Granada", "Sin especificar", "Segunda", 220, 14400, 14400, 14400, "`$/caja 18 kilos granel", "Provincia de Limarí", 800, 18),
    @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44649, 13, "Fruta", 100104, "Frutos de pepita", 100104001, "Granada", "Sin especificar", "Especial", 220, 21600, 21600, 21600, "`$/caja 18 kilos granel", "Provincia de Limarí", 1200, 18),
    @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44649, 13, "Fruta", 100104, "Frutos de pepita", 100104001, "Granada", "Sin especificar", "Primera", 250, 16200, 16200, 16200, "`$/caja 18 kilos granel", "Provincia de Limarí", 900, 18),
    @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44649, 13, "Fruta", 100104, "Frutos de pepita", 100104001, "Granada", "Sin especificar", "Segunda", 180, 14400, 14400, 14400, "`$/caja 18 kilos granel", "Provincia de Limarí", 800, 18),
    @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44644, 13, "Fruta", 100104, "Frutos de pepita", 100104001, "Granada", "Sin especificar", "Especial", 180, 18000, 18000, 18000, "`$/caja 15 kilos granel", "Provincia de Limarí", 1200, 15),
    @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44644, 13, "Fruta", 100104, "Frutos de pepita", 100104001, "Granada", "Sin especificar", "Primera", 220, 13500, 13500, 13500, "`$/caja 15 kilos granel", "Provincia de Limarí", 900, 15),
    @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44644, 13, "Fruta", 100104, "Frutos de pepita", 100104001, "Granada", "Sin especificar", "Segunda", 290, 12000, 12000, 12000, "`$/caja 15 kilos granel", "Provincia de Limarí", 800, 15),
    @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44305, 13, "Fruta", 100104, "Frutos de pepita", 100104001, "Granada", "Wonderfull", "Primera", 50, 18000, 18000, 18000, "`$/caja 15 kilos granel", "Región de O'Higgins", 1200, 15),
    @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44305, 13, "Fruta", 100104, "Frutos de pepita", 100104001, "Granada", "Wonderfull", "Segunda", 60, 15000, 15000, 15000, "`$/caja 15 kilos granel", "Región de O'Higgins", 1000, 15),
    @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44285, 13, "Fruta", 100104, "Frutos de pepita", 100104001, "Granada", "Wonderfull", "Especial", 40, 18000, 18000, 18000, "`$/caja 15 kilos empedrada", "Provincia del Elquí", 1200, 15),
    @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44285, 13, "Fruta", 100104, "Frutos de pepita", 100104001, "Granada", "Wonderfull", "Primera", 90, 15000, 15000, 15000, "`$/caja 15 kilos empedrada", "Provincia del Elquí", 1000, 15),
    @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44285, 13, "Fruta", 100104, "Frutos de pepita", 100104001, "Granada", "Wonderfull", "Segunda", 75, 12000, 12000, 12000, "`$/caja 15 kilos empedrada", "Provincia del Elquí", 800, 15),
    @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44309, 13, "Fruta", 100104, "Frutos de pepita", 100104001, "Granada", "Wonderfull", "Primera", 40, 18000, 18000, 18000, "`$/caja 15 kilos granel", "Región de O'Higgins", 1200, 15),
    @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44309, 13, "Fruta", 100104, "Frutos de pepita", 100104001, "Granada", "Wonderfull", "Segunda", 70, 15000, 15000, 15000, "`$/caja 15 kilos granel", "Región de O'Higgins", 1000, 15)
)

$arr = New-Object 'object[,]' 19,20
for ($r = 0; $r -lt 19; $r++) {
    for ($c = 0; $c -lt 20; $c++) {
        $arr[$r, $c] = $final[$r][$c]
    }
}

$ws.Range("A2:T20").Value = $arr

# Column D (Fecha) keeps the workbook's date number format on every
# data row, including the 3 newly-appended ones.
$ws.Range("D2:D20").NumberFormat = "YYYY-MM-DD HH:MM:SS"
